$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.142413774921593
$ws.Range("D2").Value = 0.02385247378005317
$ws.Range("E2").Value = 0.1636080810366707
$ws.Range("F2").Value = 0.4506741034633777
$ws.Range("G2").Value = 0.2954979547034569
$ws.Range("H2").Value = 0.4680169612718217
$ws.Range("K2").Value = 0.6343608047112923
$ws.Range("M2").Value = 0.2599284066513761
$ws.Range("O2").Value = 1.453487832759862
$ws.Range("B3").Value = 0.1329423308187501
$ws.Range("D3").Value = 0.02116874622711151
$ws.Range("E3").Value = 0.1568622348986963
$ws.Range("F3").Value = 0.4498970512007432
$ws.Range("G3").Value = 0.2959931267280993
$ws.Range("H3").Value = 0.4719883458770653
$ws.Range("K3").Value = 0.561478374576069
$ws.Range("M3").Value = 0.2290251281312408
$ws.Range("O3").Value = 1.462599218654262
$ws.Range("B4").Value = 0.1271965829031245
$ws.Range("D4").Value = 0.01951168316361418
$ws.Range("E4").Value = 0.1528795945473505
$ws.Range("F4").Value = 0.4497835145380549
$ws.Range("G4").Value = 0.2965901043679082
$ws.Range("H4").Value = 0.474687383212725
$ws.Range("K4").Value = 0.5165045148913237
$ws.Range("M4").Value = 0.2100548491284684
$ws.Range("O4").Value = 1.469351443595599
$ws.Range("B5").Value = 0.124872841543791
$ws.Range("D5").Value = 0.01883413589204963
$ws.Range("E5").Value = 0.1512964295268873
$ws.Range("F5").Value = 0.4498285748792625
$ws.Range("G5").Value = 0.2969068524993475
$ws.Range("H5").Value = 0.4758527833578228
$ws.Range("K5").Value = 0.4981223031839477
$ws.Range("M5").Value = 0.2023255689110002
$ws.Range("O5").Value = 1.472393749490394
$ws.Range("B6").Value = 0.1244880595624096
$ws.Range("D6").Value = 0.01872149322969108
$ws.Range("E6").Value = 0.1510359416015561
$ws.Range("F6").Value = 0.4498415702187089
$ws.Range("G6").Value = 0.2969638804850376
$ws.Range("H6").Value = 0.4760502542311471
$ws.Range("K6").Value = 0.4950666665654637
$ws.Range("M6").Value = 0.2010422082917245
$ws.Range("O6").Value = 1.472916468520907
$ws.Range("B7").Value = 0.1271651722501446
$ws.Range("D7").Value = 0.01950255470020323
$ws.Range("E7").Value = 0.1528580826702708
$ws.Range("F7").Value = 0.4497837525855815
$ws.Range("G7").Value = 0.2965940788962271
$ws.Range("H7").Value = 0.4747028349271076
$ws.Range("K7").Value = 0.5162568269587666
$ws.Range("M7").Value = 0.2099506040510306
$ws.Range("O7").Value = 1.469391296642712
$ws.Range("B8").Value = 0.1391336335988882
$ws.Range("D8").Value = 0.02292906842919251
$ws.Range("E8").Value = 0.1612488763230573
$ws.Range("F8").Value = 0.450330655739144
$ws.Range("G8").Value = 0.2956077770030774
$ws.Range("H8").Value = 0.4693322045099535
$ws.Range("K8").Value = 0.6092781179497422
$ws.Range("M8").Value = 0.249272105722298
$ws.Range("O8").Value = 1.456388884604195
$ws.Range("B9").Value = 0.1631516676451525
$ws.Range("D9").Value = 0.02957349319558489
$ws.Range("E9").Value = 0.1789807880906622
$ws.Range("F9").Value = 0.4542931583769558
$ws.Range("G9").Value = 0.2960068771553352
$ws.Range("H9").Value = 0.4608685534058097
$ws.Range("K9").Value = 0.7898702725357225
$ws.Range("M9").Value = 0.3264145258855464
$ws.Range("O9").Value = 1.440097882490633
$ws.Range("B10").Value = 0.1811261772330255
$ws.Range("D10").Value = 0.03440773518208573
$ws.Range("E10").Value = 0.1928071712868942
$ws.Range("F10").Value = 0.4589746739484255
$ws.Range("G10").Value = 0.2977354345449754
$ws.Range("H10").Value = 0.455911947941118
$ws.Range("K10").Value = 0.9213900797547581
$ws.Range("M10").Value = 0.383113896132258
$ws.Range("O10").Value = 1.433771102881678
$ws.Range("B11").Value = 0.1893735033360855
$ws.Range("D11").Value = 0.0365963242264371
$ws.Range("E11").Value = 0.1992749492471191
$ws.Range("F11").Value = 0.4614906803834415
$ws.Range("G11").Value = 0.2988362710111545
$ws.Range("H11").Value = 0.4539312144993488
$ws.Range("K11").Value = 0.9809599195114629
$ws.Range("M11").Value = 0.408913805830224
$ws.Range("O11").Value = 1.432124482465213
$ws.Range("B12").Value = 0.1925065745647601
$ws.Range("D12").Value = 0.03742353302348533
$ws.Range("E12").Value = 0.2017500564312229
$ws.Range("F12").Value = 0.4624991088094816
$ws.Range("G12").Value = 0.2992985661507674
$ws.Range("H12").Value = 0.4532205865493921
$ws.Range("K12").Value = 1.003479108586987
$ws.Range("M12").Value = 0.4186845248072188
$ws.Range("O12").Value = 1.431678520680805
$ws.Range("B13").Value = 0.1918313700975744
$ws.Range("D13").Value = 0.03724544905718119
$ws.Range("E13").Value = 0.2012158412659701
$ws.Range("F13").Value = 0.4622794479227039
$ws.Range("G13").Value = 0.299196978737001
$ws.Range("H13").Value = 0.4533718788180892
$ws.Range("K13").Value = 0.9986309371127504
$ws.Range("M13").Value = 0.4165801903334625
$ws.Range("O13").Value = 1.431766660879418
$ws.Range("B14").Value = 0.1896310637885961
$ws.Range("D14").Value = 0.0366644107789682
$ws.Range("E14").Value = 0.1994780570612917
$ws.Range("F14").Value = 0.4615725280178395
$ws.Range("G14").Value = 0.2988733926028715
$ws.Range("H14").Value = 0.4538719602774677
$ws.Range("K14").Value = 0.9828133683587339
$ws.Range("M14").Value = 0.4097176326063021
$ws.Range("O14").Value = 1.432084231180141
$ws.Range("B15").Value = 0.1882846083351239
$ws.Range("D15").Value = 0.0363083031693634
$ws.Range("E15").Value = 0.1984169956511508
$ws.Range("F15").Value = 0.4611467725491352
$ws.Range("G15").Value = 0.2986811094768314
$ws.Range("H15").Value = 0.4541834106548492
$ws.Range("K15").Value = 0.9731195819778975
$ws.Range("M15").Value = 0.4055142240584928
$ws.Range("O15").Value = 1.432301892627464
$ws.Range("B16").Value = 0.1805886019286191
$ws.Range("D16").Value = 0.03426448971472951
$ws.Range("E16").Value = 0.1923880952627925
$ws.Range("F16").Value = 0.4588180312984989
$ws.Range("G16").Value = 0.297669838669961
$ws.Range("H16").Value = 0.4560469102457603
$ws.Range("K16").Value = 0.9174917199286767
$ws.Range("M16").Value = 0.3814279397130349
$ws.Range("O16").Value = 1.433903535962173
$ws.Range("B17").Value = 0.1758853215743841
$ws.Range("D17").Value = 0.0330079447656999
$ws.Range("E17").Value = 0.1887353798733074
$ws.Range("F17").Value = 0.4574884614358723
$ws.Range("G17").Value = 0.2971301570362996
$ws.Range("H17").Value = 0.4572603155211254
$ws.Range("K17").Value = 0.8832984848514513
$ws.Range("M17").Value = 0.3666534803454127
$ws.Range("O17").Value = 1.435201852744171
$ws.Range("B18").Value = 0.1731867773933686
$ws.Range("D18").Value = 0.03228422488979987
$ws.Range("E18").Value = 0.1866511795778365
$ws.Range("F18").Value = 0.4567600875412907
$ws.Range("G18").Value = 0.2968493370779655
$ws.Range("H18").Value = 0.4579840276790392
$ws.Range("K18").Value = 0.8636071250188024
$ws.Range("M18").Value = 0.3581562784783117
$ws.Range("O18").Value = 1.436064492030766
$ws.Range("B19").Value = 0.1722742454109749
$ws.Range("D19").Value = 0.03203901746900328
$ws.Range("E19").Value = 0.1859483722537121
$ws.Range("F19").Value = 0.4565197139660668
$ws.Range("G19").Value = 0.296759332168925
$ws.Range("H19").Value = 0.4582334934956123
$ws.Range("K19").Value = 0.8569358339227335
$ws.Range("M19").Value = 0.3552793931488196
$ws.Range("O19").Value = 1.436376454554633
$ws.Range("B20").Value = 0.1763853059132146
$ws.Range("D20").Value = 0.03314180876905937
$ws.Range("E20").Value = 0.1891224825086368
$ws.Range("F20").Value = 0.4576262325181304
$ws.Range("G20").Value = 0.297184543045617
$ws.Range("H20").Value = 0.4571284767446215
$ws.Range("K20").Value = 0.8869409394354193
$ws.Range("M20").Value = 0.3682261779651341
$ws.Range("O20").Value = 1.435051647775907
$ws.Range("B21").Value = 0.1902770773674831
$ws.Range("D21").Value = 0.03683511858936583
$ws.Range("E21").Value = 0.1999877811663922
$ws.Range("F21").Value = 0.4617786557316919
$ws.Range("G21").Value = 0.2989672030377193
$ws.Range("H21").Value = 0.4537240038404633
$ws.Range("K21").Value = 0.9874604304474701
$ws.Range("M21").Value = 0.4117333103206846
$ws.Range("O21").Value = 1.431986129754677
$ws.Range("B22").Value = 0.1994142912471659
$ws.Range("D22").Value = 0.03923978119090066
$ws.Range("E22").Value = 0.2072399822282236
$ws.Range("F22").Value = 0.4648170510699643
$ws.Range("G22").Value = 0.3003971831087711
$ws.Range("H22").Value = 0.4517288399317465
$ws.Range("K22").Value = 1.052930024609566
$ws.Range("M22").Value = 0.4401726329444813
$ws.Range("O22").Value = 1.43101792181352
$ws.Range("B23").Value = 0.1945323244411981
$ws.Range("D23").Value = 0.03795721859221146
$ws.Range("E23").Value = 0.2033554241442985
$ws.Range("F23").Value = 0.4631656693395954
$ws.Range("G23").Value = 0.2996096688365526
$ws.Range("H23").Value = 0.4527726565408017
$ws.Range("K23").Value = 1.018008764432409
$ws.Range("M23").Value = 0.4249936375358203
$ws.Range("O23").Value = 1.431439781443913
$ws.Range("B24").Value = 0.1761592460541976
$ws.Range("D24").Value = 0.03308129294718753
$ws.Range("E24").Value = 0.1889474242329356
$ws.Range("F24").Value = 0.4575638340409185
$ws.Range("G24").Value = 0.2971598634109327
$ws.Range("H24").Value = 0.4571879997276653
$ws.Range("K24").Value = 0.8852942890742668
$ws.Range("M24").Value = 0.3675151712317657
$ws.Range("O24").Value = 1.435119193461077
$ws.Range("B25").Value = 0.1565960697809459
$ws.Range("D25").Value = 0.02778420341057597
$ws.Range("E25").Value = 0.1740450039608987
$ws.Range("F25").Value = 0.4529109221252838
$ws.Range("G25").Value = 0.2956477786064156
$ws.Range("H25").Value = 0.4629366990435599
$ws.Range("K25").Value = 0.7412156407458497
$ws.Range("M25").Value = 0.305541644668935
$ws.Range("O25").Value = 1.443516321584724
